$d = $word.ActiveDocument

$pairs = @(
    @("720×7=5040", "460×7=3220"),
    @("617×8=4936", "746×4=2984"),
    @("920×2=1840", "743×3=2229"),
    @("821×5=4105", "867×3=2601"),
    @("233×4=932", "753×9=6777"),
    @("396×9=3564", "746×7=5222"),
    @("252×2=504", "909×4=3636"),
    @("305×2=610", "845×5=4225"),
    @("709×4=2836", "681×5=3405"),
    @("186×9=1674", "614×2=1228"),
    @("214×7=1498", "649×3=1947"),
    @("921×8=7368", "389×7=2723"),
    @("115×4=460", "265×5=1325"),
    @("624×9=5616", "673×2=1346"),
    @("578×7=4046", "664×5=3320"),
    @("738×3=2214", "563×9=5067"),
    @("328×5=1640", "916×6=5496"),
    @("723×5=3615", "374×8=2992"),
    @("620×6=3720", "479×4=1916"),
    @("988×8=7904", "973×7=6811"),
    @("650×3=1950", "495×7=3465"),
    @("949×8=7592", "411×8=3288"),
    @("308×9=2772", "811×6=4866"),
    @("555×8=4440", "159×8=1272"),
    @("708×8=5664", "851×7=5957")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
